# "Frogger Language Definition" doc -- add the assignment operator (=) to
# the OPERATORS paragraph's list of arithmetic operators.
#
# Before: "...multiplication (**), and division (//). Boolean operators: ..."
# After:  "...multiplication (**), division (//), and assignment (=). Boolean operators: ..."

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "and division (//). ",   # old text
    $false,                  # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "division (//), and assignment (=). ",  # new text
    2)                       # Replace (wdReplaceAll)

if (-not $found) {
    throw "Could not locate the OPERATORS text to update."
}
